$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values are not auto-converted to numbers by Excel
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.281.70'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.53'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.77'
$ws.Range("E5").Value = '  -1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.596'
$ws.Range("E6").Value = '  +2.96%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.16'
$ws.Range("E8").Value = '  +3.22%  '
$ws.Range("E9").Value = '  -4.21%  '
$ws.Range("E10").Value = '  -4.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0962'
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.050.20'
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.10'
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.821.01'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.626'
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.271.89'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.50'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.70'
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0767'
$ws.Range("E20").Value = '  -4.33%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.21'
$ws.Range("E21").Value = '  -3.46%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.44'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.97'
$ws.Range("E26").Value = '  +5.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.14'
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.119'
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.76'
$ws.Range("E31").Value = '  -2.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.86'
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0512'
$ws.Range("E33").Value = '  -3.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.75'
$ws.Range("E34").Value = '  -5.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.356.55'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.642'
$ws.Range("E36").Value = '  -5.56%  '
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.34'
$ws.Range("E38").Value = '  -8.45%  '
$ws.Range("E39").Value = '  -4.06%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '80.23'
$ws.Range("E42").Value = '  -3.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.929'
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.16'
$ws.Range("E44").Value = '  +5.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.12'
$ws.Range("E45").Value = '  -5.25%  '
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.951.93'
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.74'
$ws.Range("E48").Value = '  -5.11%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.50'
$ws.Range("E50").Value = '  -3.49%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0120'
$ws.Range("E51").Value = '  -8.31%  '
